# add the NA's under duplicate_image_filename
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("duplicate_image_filename") gets "NA" for data rows 2 through 21.
$ws.Range("E2:E21").Value = "NA"
